$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-72 (only changed cells) ---

# Row 2
$ws.Cells.Item(2, 4).Value = 44225
$ws.Cells.Item(2, 9).Value = '1a (nueva lavada)'
$ws.Cells.Item(2, 11).Value = 12000
$ws.Cells.Item(2, 12).Value = 13000
$ws.Cells.Item(2, 13).Value = 12500
$ws.Cells.Item(2, 15).Value = 'Región del Maule'
$ws.Cells.Item(2, 16).Value = 500

# Row 3
$ws.Cells.Item(3, 4).Value = 44225
$ws.Cells.Item(3, 8).Value = 'Rosara'
$ws.Cells.Item(3, 9).Value = '1a (cosecha)'
$ws.Cells.Item(3, 11).Value = 9000
$ws.Cells.Item(3, 12).Value = 10000
$ws.Cells.Item(3, 13).Value = 9500
$ws.Cells.Item(3, 15).Value = 'Región del Maule'
$ws.Cells.Item(3, 16).Value = 380

# Row 4
$ws.Cells.Item(4, 4).Value = 44215
$ws.Cells.Item(4, 8).Value = 'Pukará'
$ws.Cells.Item(4, 9).Value = '1a nueva(o)'
$ws.Cells.Item(4, 11).Value = 11000
$ws.Cells.Item(4, 12).Value = 12000
$ws.Cells.Item(4, 13).Value = 11500
$ws.Cells.Item(4, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(4, 16).Value = 460

# Row 5
$ws.Cells.Item(5, 4).Value = 44592
$ws.Cells.Item(5, 8).Value = 'Patagonia'
$ws.Cells.Item(5, 9).Value = '1a (cosecha)'
$ws.Cells.Item(5, 11).Value = 10000
$ws.Cells.Item(5, 12).Value = 11000
$ws.Cells.Item(5, 13).Value = 10500
$ws.Cells.Item(5, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(5, 15).Value = 'Región del Maule'
$ws.Cells.Item(5, 16).Value = 420

# Row 6
$ws.Cells.Item(6, 4).Value = 44257
$ws.Cells.Item(6, 8).Value = 'Asterix'
$ws.Cells.Item(6, 9).Value = '1a (nueva lavada)'
$ws.Cells.Item(6, 11).Value = 9000
$ws.Cells.Item(6, 12).Value = 9500
$ws.Cells.Item(6, 13).Value = 9250
$ws.Cells.Item(6, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(6, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(6, 16).Value = 370

# Row 7
$ws.Cells.Item(7, 4).Value = 44525
$ws.Cells.Item(7, 8).Value = 'Asterix'
$ws.Cells.Item(7, 9).Value = '1a (nueva lavada)'
$ws.Cells.Item(7, 11).Value = 15000
$ws.Cells.Item(7, 12).Value = 16000
$ws.Cells.Item(7, 13).Value = 15500
$ws.Cells.Item(7, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(7, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(7, 16).Value = 620

# Row 8
$ws.Cells.Item(8, 4).Value = 44615
$ws.Cells.Item(8, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(8, 11).Value = 11000
$ws.Cells.Item(8, 12).Value = 12000
$ws.Cells.Item(8, 13).Value = 11500
$ws.Cells.Item(8, 16).Value = 460

# Row 9
$ws.Cells.Item(9, 4).Value = 44400
$ws.Cells.Item(9, 8).Value = 'Cardinal'
$ws.Cells.Item(9, 9).Value = '1a nueva(o)'
$ws.Cells.Item(9, 11).Value = 12000
$ws.Cells.Item(9, 12).Value = 13000
$ws.Cells.Item(9, 13).Value = 12500
$ws.Cells.Item(9, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(9, 16).Value = 500

# Row 10
$ws.Cells.Item(10, 4).Value = 44357
$ws.Cells.Item(10, 8).Value = 'Asterix'
$ws.Cells.Item(10, 9).Value = '1a (nueva lavada)'
$ws.Cells.Item(10, 11).Value = 9000
$ws.Cells.Item(10, 12).Value = 9500
$ws.Cells.Item(10, 13).Value = 9250
$ws.Cells.Item(10, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(10, 16).Value = 370

# Row 11
$ws.Cells.Item(11, 4).Value = 44357
$ws.Cells.Item(11, 8).Value = 'Cardinal'
$ws.Cells.Item(11, 11).Value = 10000
$ws.Cells.Item(11, 12).Value = 11000
$ws.Cells.Item(11, 13).Value = 10500
$ws.Cells.Item(11, 16).Value = 420

# Row 12
$ws.Cells.Item(12, 4).Value = 44644
$ws.Cells.Item(12, 8).Value = 'Asterix'
$ws.Cells.Item(12, 9).Value = '1a (cosecha)'
$ws.Cells.Item(12, 11).Value = 8500
$ws.Cells.Item(12, 12).Value = 9000
$ws.Cells.Item(12, 13).Value = 8750
$ws.Cells.Item(12, 16).Value = 350

# Row 13
$ws.Cells.Item(13, 4).Value = 44299
$ws.Cells.Item(13, 9).Value = '1a (cosecha)'
$ws.Cells.Item(13, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(13, 15).Value = 'Región de Los Lagos'

# Row 14
$ws.Cells.Item(14, 4).Value = 44433
$ws.Cells.Item(14, 9).Value = '1a (guarda)'
$ws.Cells.Item(14, 11).Value = 9000
$ws.Cells.Item(14, 12).Value = 9500
$ws.Cells.Item(14, 13).Value = 9250
$ws.Cells.Item(14, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(14, 16).Value = 370

# Row 15
$ws.Cells.Item(15, 4).Value = 44614
$ws.Cells.Item(15, 9).Value = '1a (cosecha)'
$ws.Cells.Item(15, 11).Value = 9000
$ws.Cells.Item(15, 12).Value = 10000
$ws.Cells.Item(15, 13).Value = 9500
$ws.Cells.Item(15, 15).Value = 'Región del Maule'
$ws.Cells.Item(15, 16).Value = 380

# Row 16
$ws.Cells.Item(16, 4).Value = 44497
$ws.Cells.Item(16, 8).Value = 'Asterix'
$ws.Cells.Item(16, 9).Value = '1a (cosecha)'
$ws.Cells.Item(16, 11).Value = 15000
$ws.Cells.Item(16, 12).Value = 16000
$ws.Cells.Item(16, 13).Value = 15500
$ws.Cells.Item(16, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(16, 16).Value = 620

# Row 17
$ws.Cells.Item(17, 4).Value = 44322
$ws.Cells.Item(17, 8).Value = 'Rodeo'
$ws.Cells.Item(17, 9).Value = '1a (cosecha)'
$ws.Cells.Item(17, 11).Value = 8000
$ws.Cells.Item(17, 12).Value = 8500
$ws.Cells.Item(17, 13).Value = 8250
$ws.Cells.Item(17, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(17, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(17, 16).Value = 330

# Row 18
$ws.Cells.Item(18, 4).Value = 44427
$ws.Cells.Item(18, 8).Value = 'Asterix'
$ws.Cells.Item(18, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(18, 11).Value = 9000
$ws.Cells.Item(18, 12).Value = 10000
$ws.Cells.Item(18, 13).Value = 9500
$ws.Cells.Item(18, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(18, 16).Value = 380

# Row 19
$ws.Cells.Item(19, 4).Value = 44427
$ws.Cells.Item(19, 8).Value = 'Cardinal'
$ws.Cells.Item(19, 9).Value = '1a (cosecha)'
$ws.Cells.Item(19, 15).Value = 'Provincia de Melipilla'

# Row 20
$ws.Cells.Item(20, 4).Value = 44692
$ws.Cells.Item(20, 9).Value = '1a (cosecha)'
$ws.Cells.Item(20, 11).Value = 8500
$ws.Cells.Item(20, 12).Value = 9000
$ws.Cells.Item(20, 13).Value = 8750
$ws.Cells.Item(20, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(20, 16).Value = 350

# Row 21
$ws.Cells.Item(21, 4).Value = 44574
$ws.Cells.Item(21, 9).Value = '1a (cosecha)'
$ws.Cells.Item(21, 11).Value = 13000
$ws.Cells.Item(21, 12).Value = 14000
$ws.Cells.Item(21, 13).Value = 13500
$ws.Cells.Item(21, 15).Value = 'Región del Maule'
$ws.Cells.Item(21, 16).Value = 540

# Row 22
$ws.Cells.Item(22, 4).Value = 44617
$ws.Cells.Item(22, 8).Value = 'Asterix'
$ws.Cells.Item(22, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(22, 11).Value = 11000
$ws.Cells.Item(22, 12).Value = 12000
$ws.Cells.Item(22, 13).Value = 11500
$ws.Cells.Item(22, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(22, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(22, 16).Value = 460

# Row 23
$ws.Cells.Item(23, 4).Value = 44617
$ws.Cells.Item(23, 8).Value = 'Rosara'
$ws.Cells.Item(23, 11).Value = 9000
$ws.Cells.Item(23, 12).Value = 10000
$ws.Cells.Item(23, 13).Value = 9500
$ws.Cells.Item(23, 16).Value = 380

# Row 24
$ws.Cells.Item(24, 4).Value = 44411
$ws.Cells.Item(24, 11).Value = 12000
$ws.Cells.Item(24, 12).Value = 13000
$ws.Cells.Item(24, 13).Value = 12500
$ws.Cells.Item(24, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(24, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(24, 16).Value = 500

# Row 25
$ws.Cells.Item(25, 4).Value = 44159
$ws.Cells.Item(25, 8).Value = 'Rosara'
$ws.Cells.Item(25, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(25, 11).Value = 11000
$ws.Cells.Item(25, 12).Value = 12000
$ws.Cells.Item(25, 13).Value = 11500
$ws.Cells.Item(25, 16).Value = 460

# Row 26
$ws.Cells.Item(26, 4).Value = 44559
$ws.Cells.Item(26, 9).Value = '1a nueva(o)'
$ws.Cells.Item(26, 11).Value = 13000
$ws.Cells.Item(26, 12).Value = 14000
$ws.Cells.Item(26, 13).Value = 13500
$ws.Cells.Item(26, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(26, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(26, 16).Value = 540

# Row 27
$ws.Cells.Item(27, 4).Value = 44181
$ws.Cells.Item(27, 9).Value = '1a nueva(o)'
$ws.Cells.Item(27, 11).Value = 13000
$ws.Cells.Item(27, 12).Value = 14000
$ws.Cells.Item(27, 13).Value = 13500
$ws.Cells.Item(27, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(27, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(27, 16).Value = 540

# Row 28
$ws.Cells.Item(28, 4).Value = 44218
$ws.Cells.Item(28, 8).Value = 'Asterix'
$ws.Cells.Item(28, 9).Value = '1a (nueva lavada)'
$ws.Cells.Item(28, 10).Value = 1000
$ws.Cells.Item(28, 11).Value = 13000
$ws.Cells.Item(28, 12).Value = 14000
$ws.Cells.Item(28, 13).Value = 13500
$ws.Cells.Item(28, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(28, 16).Value = 540

# Row 29
$ws.Cells.Item(29, 4).Value = 44651
$ws.Cells.Item(29, 11).Value = 9000
$ws.Cells.Item(29, 12).Value = 10000
$ws.Cells.Item(29, 13).Value = 9500
$ws.Cells.Item(29, 16).Value = 380

# Row 30
$ws.Cells.Item(30, 4).Value = 44434
$ws.Cells.Item(30, 9).Value = '1a (guarda)'
$ws.Cells.Item(30, 11).Value = 9000
$ws.Cells.Item(30, 12).Value = 9500
$ws.Cells.Item(30, 13).Value = 9250
$ws.Cells.Item(30, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(30, 16).Value = 370

# Row 31
$ws.Cells.Item(31, 4).Value = 44469
$ws.Cells.Item(31, 8).Value = 'Asterix'
$ws.Cells.Item(31, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(31, 11).Value = 10000
$ws.Cells.Item(31, 12).Value = 11000
$ws.Cells.Item(31, 13).Value = 10500
$ws.Cells.Item(31, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(31, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(31, 16).Value = 420

# Row 32
$ws.Cells.Item(32, 4).Value = 44707
$ws.Cells.Item(32, 11).Value = 9500
$ws.Cells.Item(32, 12).Value = 10000
$ws.Cells.Item(32, 13).Value = 9750
$ws.Cells.Item(32, 16).Value = 390

# Row 33
$ws.Cells.Item(33, 4).Value = 44384
$ws.Cells.Item(33, 8).Value = 'Rodeo'
$ws.Cells.Item(33, 9).Value = '1a (guarda)'
$ws.Cells.Item(33, 10).Value = 600
$ws.Cells.Item(33, 11).Value = 8000
$ws.Cells.Item(33, 12).Value = 8500
$ws.Cells.Item(33, 13).Value = 8208
$ws.Cells.Item(33, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(33, 16).Value = 328

# Row 34
$ws.Cells.Item(34, 4).Value = 44319
$ws.Cells.Item(34, 8).Value = 'Asterix'
$ws.Cells.Item(34, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(34, 11).Value = 9000
$ws.Cells.Item(34, 12).Value = 9500
$ws.Cells.Item(34, 13).Value = 9250
$ws.Cells.Item(34, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(34, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(34, 16).Value = 370

# Row 35
$ws.Cells.Item(35, 4).Value = 44386
$ws.Cells.Item(35, 12).Value = 9500
$ws.Cells.Item(35, 13).Value = 9250
$ws.Cells.Item(35, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(35, 16).Value = 370

# Row 36
$ws.Cells.Item(36, 4).Value = 44208
$ws.Cells.Item(36, 8).Value = 'Asterix'
$ws.Cells.Item(36, 9).Value = '1a nueva(o)'
$ws.Cells.Item(36, 11).Value = 16000
$ws.Cells.Item(36, 12).Value = 17000
$ws.Cells.Item(36, 13).Value = 16500
$ws.Cells.Item(36, 15).Value = 'Región del Maule'
$ws.Cells.Item(36, 16).Value = 660

# Row 37
$ws.Cells.Item(37, 4).Value = 44301
$ws.Cells.Item(37, 8).Value = 'Asterix'
$ws.Cells.Item(37, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(37, 11).Value = 8500
$ws.Cells.Item(37, 12).Value = 9000
$ws.Cells.Item(37, 13).Value = 8750
$ws.Cells.Item(37, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(37, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(37, 16).Value = 350

# Row 38
$ws.Cells.Item(38, 4).Value = 44309
$ws.Cells.Item(38, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(38, 11).Value = 9000
$ws.Cells.Item(38, 12).Value = 9500
$ws.Cells.Item(38, 13).Value = 9250
$ws.Cells.Item(38, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(38, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(38, 16).Value = 370

# Row 39
$ws.Cells.Item(39, 4).Value = 44166
$ws.Cells.Item(39, 8).Value = 'Cardinal'
$ws.Cells.Item(39, 9).Value = '1a nueva(o)'
$ws.Cells.Item(39, 11).Value = 13000
$ws.Cells.Item(39, 12).Value = 14000
$ws.Cells.Item(39, 13).Value = 13500
$ws.Cells.Item(39, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(39, 16).Value = 540

# Row 40
$ws.Cells.Item(40, 4).Value = 44571
$ws.Cells.Item(40, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(40, 11).Value = 14000
$ws.Cells.Item(40, 12).Value = 15000
$ws.Cells.Item(40, 13).Value = 14500
$ws.Cells.Item(40, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(40, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(40, 16).Value = 580

# Row 41
$ws.Cells.Item(41, 4).Value = 44238
$ws.Cells.Item(41, 8).Value = 'Patagonia'
$ws.Cells.Item(41, 9).Value = '1a nueva(o)'
$ws.Cells.Item(41, 11).Value = 8500
$ws.Cells.Item(41, 12).Value = 9000
$ws.Cells.Item(41, 13).Value = 8750
$ws.Cells.Item(41, 16).Value = 350

# Row 42
$ws.Cells.Item(42, 4).Value = 44222
$ws.Cells.Item(42, 9).Value = '1a nueva(o)'
$ws.Cells.Item(42, 11).Value = 10000
$ws.Cells.Item(42, 12).Value = 11000
$ws.Cells.Item(42, 13).Value = 10500
$ws.Cells.Item(42, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(42, 16).Value = 420

# Row 43
$ws.Cells.Item(43, 4).Value = 44677
$ws.Cells.Item(43, 8).Value = 'Asterix'
$ws.Cells.Item(43, 9).Value = '1a (cosecha)'
$ws.Cells.Item(43, 10).Value = 1000
$ws.Cells.Item(43, 11).Value = 8500
$ws.Cells.Item(43, 12).Value = 9000
$ws.Cells.Item(43, 13).Value = 8750
$ws.Cells.Item(43, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(43, 16).Value = 350

# Row 44
$ws.Cells.Item(44, 4).Value = 44336
$ws.Cells.Item(44, 8).Value = 'Cardinal'
$ws.Cells.Item(44, 9).Value = '1a nueva(o)'
$ws.Cells.Item(44, 10).Value = 500
$ws.Cells.Item(44, 11).Value = 11000
$ws.Cells.Item(44, 12).Value = 12000
$ws.Cells.Item(44, 13).Value = 11500
$ws.Cells.Item(44, 15).Value = 'La Ligua'
$ws.Cells.Item(44, 16).Value = 460

# Row 45
$ws.Cells.Item(45, 4).Value = 44334
$ws.Cells.Item(45, 8).Value = 'Asterix'
$ws.Cells.Item(45, 12).Value = 9000
$ws.Cells.Item(45, 13).Value = 8500
$ws.Cells.Item(45, 16).Value = 340

# Row 46
$ws.Cells.Item(46, 4).Value = 44232
$ws.Cells.Item(46, 8).Value = 'Asterix'
$ws.Cells.Item(46, 9).Value = '1a (nueva lavada)'
$ws.Cells.Item(46, 11).Value = 11000
$ws.Cells.Item(46, 12).Value = 12000
$ws.Cells.Item(46, 13).Value = 11500
$ws.Cells.Item(46, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(46, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(46, 16).Value = 460

# Row 47
$ws.Cells.Item(47, 4).Value = 44251
$ws.Cells.Item(47, 8).Value = 'Patagonia'
$ws.Cells.Item(47, 9).Value = '1a nueva(o)'
$ws.Cells.Item(47, 11).Value = 8500
$ws.Cells.Item(47, 12).Value = 9000
$ws.Cells.Item(47, 13).Value = 8750
$ws.Cells.Item(47, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(47, 16).Value = 350

# Row 49
$ws.Cells.Item(49, 4).Value = 44551
$ws.Cells.Item(49, 9).Value = '1a (cosecha)'
$ws.Cells.Item(49, 11).Value = 15000
$ws.Cells.Item(49, 12).Value = 16000
$ws.Cells.Item(49, 13).Value = 15500
$ws.Cells.Item(49, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(49, 15).Value = 'Región del Maule'
$ws.Cells.Item(49, 16).Value = 620

# Row 50
$ws.Cells.Item(50, 4).Value = 44498
$ws.Cells.Item(50, 8).Value = 'Rosara'
$ws.Cells.Item(50, 9).Value = '1a nueva(o)'
$ws.Cells.Item(50, 11).Value = 14000
$ws.Cells.Item(50, 12).Value = 15000
$ws.Cells.Item(50, 13).Value = 14500
$ws.Cells.Item(50, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(50, 15).Value = 'Región del Maule'
$ws.Cells.Item(50, 16).Value = 580

# Row 51
$ws.Cells.Item(51, 4).Value = 44545
$ws.Cells.Item(51, 9).Value = '1a (cosecha)'
$ws.Cells.Item(51, 11).Value = 14000
$ws.Cells.Item(51, 12).Value = 15000
$ws.Cells.Item(51, 13).Value = 14500
$ws.Cells.Item(51, 15).Value = 'Región del Maule'
$ws.Cells.Item(51, 16).Value = 580

# Row 52
$ws.Cells.Item(52, 4).Value = 44264
$ws.Cells.Item(52, 8).Value = 'Patagonia'
$ws.Cells.Item(52, 9).Value = '1a (cosecha)'
$ws.Cells.Item(52, 11).Value = 8000
$ws.Cells.Item(52, 12).Value = 8500
$ws.Cells.Item(52, 13).Value = 8250
$ws.Cells.Item(52, 15).Value = 'Región del Maule'
$ws.Cells.Item(52, 16).Value = 330

# Row 53
$ws.Cells.Item(53, 4).Value = 44657
$ws.Cells.Item(53, 9).Value = '1a (cosecha)'
$ws.Cells.Item(53, 11).Value = 9500
$ws.Cells.Item(53, 12).Value = 10000
$ws.Cells.Item(53, 13).Value = 9750
$ws.Cells.Item(53, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(53, 16).Value = 390

# Row 54
$ws.Cells.Item(54, 4).Value = 44235
$ws.Cells.Item(54, 8).Value = 'Rosara'
$ws.Cells.Item(54, 9).Value = '1a (cosecha)'
$ws.Cells.Item(54, 11).Value = 8000
$ws.Cells.Item(54, 12).Value = 9000
$ws.Cells.Item(54, 13).Value = 8500
$ws.Cells.Item(54, 15).Value = 'Región del Maule'
$ws.Cells.Item(54, 16).Value = 340

# Row 55
$ws.Cells.Item(55, 4).Value = 44413
$ws.Cells.Item(55, 8).Value = 'Asterix'
$ws.Cells.Item(55, 9).Value = '1a (nueva lavada)'
$ws.Cells.Item(55, 11).Value = 11000
$ws.Cells.Item(55, 12).Value = 12000
$ws.Cells.Item(55, 13).Value = 11500
$ws.Cells.Item(55, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(55, 16).Value = 460

# Row 56
$ws.Cells.Item(56, 4).Value = 44727
$ws.Cells.Item(56, 9).Value = '1a (guarda)'
$ws.Cells.Item(56, 11).Value = 9000
$ws.Cells.Item(56, 12).Value = 10000
$ws.Cells.Item(56, 13).Value = 9500
$ws.Cells.Item(56, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(56, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(56, 16).Value = 380

# Row 57
$ws.Cells.Item(57, 4).Value = 44727
$ws.Cells.Item(57, 8).Value = 'Rodeo'
$ws.Cells.Item(57, 9).Value = '1a (guarda)'
$ws.Cells.Item(57, 11).Value = 9000
$ws.Cells.Item(57, 12).Value = 9500
$ws.Cells.Item(57, 13).Value = 9250
$ws.Cells.Item(57, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(57, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(57, 16).Value = 370

# Row 58
$ws.Cells.Item(58, 4).Value = 44670
$ws.Cells.Item(58, 11).Value = 8500
$ws.Cells.Item(58, 12).Value = 9000
$ws.Cells.Item(58, 13).Value = 8750
$ws.Cells.Item(58, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(58, 16).Value = 350

# Row 59
$ws.Cells.Item(59, 4).Value = 44601
$ws.Cells.Item(59, 8).Value = 'Asterix'
$ws.Cells.Item(59, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(59, 11).Value = 11000
$ws.Cells.Item(59, 12).Value = 12000
$ws.Cells.Item(59, 13).Value = 11500
$ws.Cells.Item(59, 16).Value = 460

# Row 60
$ws.Cells.Item(60, 4).Value = 44601
$ws.Cells.Item(60, 8).Value = 'Patagonia'
$ws.Cells.Item(60, 11).Value = 9000
$ws.Cells.Item(60, 12).Value = 10000
$ws.Cells.Item(60, 13).Value = 9500
$ws.Cells.Item(60, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(60, 16).Value = 380

# Row 61
$ws.Cells.Item(61, 4).Value = 44211
$ws.Cells.Item(61, 9).Value = '1a (nueva lavada)'
$ws.Cells.Item(61, 11).Value = 15000
$ws.Cells.Item(61, 12).Value = 16000
$ws.Cells.Item(61, 13).Value = 15500
$ws.Cells.Item(61, 15).Value = 'Región del Maule'
$ws.Cells.Item(61, 16).Value = 620

# Row 62
$ws.Cells.Item(62, 4).Value = 44435
$ws.Cells.Item(62, 9).Value = '1a (guarda)'
$ws.Cells.Item(62, 11).Value = 9000
$ws.Cells.Item(62, 12).Value = 9500
$ws.Cells.Item(62, 13).Value = 9250
$ws.Cells.Item(62, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(62, 16).Value = 370

# Row 63
$ws.Cells.Item(63, 4).Value = 44435
$ws.Cells.Item(63, 8).Value = 'Asterix'
$ws.Cells.Item(63, 9).Value = '1a (guarda)'
$ws.Cells.Item(63, 11).Value = 9000
$ws.Cells.Item(63, 12).Value = 9500
$ws.Cells.Item(63, 13).Value = 9250
$ws.Cells.Item(63, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(63, 16).Value = 370

# Row 64
$ws.Cells.Item(64, 4).Value = 44216
$ws.Cells.Item(64, 8).Value = 'Rosara'
$ws.Cells.Item(64, 9).Value = '1a nueva(o)'
$ws.Cells.Item(64, 11).Value = 12000
$ws.Cells.Item(64, 12).Value = 13000
$ws.Cells.Item(64, 13).Value = 12500
$ws.Cells.Item(64, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(64, 15).Value = 'Región del Maule'
$ws.Cells.Item(64, 16).Value = 500

# Row 65
$ws.Cells.Item(65, 4).Value = 44476
$ws.Cells.Item(65, 8).Value = 'Rosara'
$ws.Cells.Item(65, 9).Value = '1a (guarda)'
$ws.Cells.Item(65, 11).Value = 11000
$ws.Cells.Item(65, 12).Value = 12000
$ws.Cells.Item(65, 13).Value = 11500
$ws.Cells.Item(65, 15).Value = 'Región del Maule'
$ws.Cells.Item(65, 16).Value = 460

# Row 66
$ws.Cells.Item(66, 4).Value = 44224
$ws.Cells.Item(66, 9).Value = '1a nueva(o)'
$ws.Cells.Item(66, 11).Value = 10000
$ws.Cells.Item(66, 12).Value = 11000
$ws.Cells.Item(66, 13).Value = 10500
$ws.Cells.Item(66, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(66, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(66, 16).Value = 420

# Row 67
$ws.Cells.Item(67, 4).Value = 44224
$ws.Cells.Item(67, 8).Value = 'Rodeo'
$ws.Cells.Item(67, 9).Value = '1a nueva(o)'

# Row 68
$ws.Cells.Item(68, 4).Value = 44449
$ws.Cells.Item(68, 8).Value = 'Asterix'
$ws.Cells.Item(68, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(68, 11).Value = 10000
$ws.Cells.Item(68, 12).Value = 11000
$ws.Cells.Item(68, 13).Value = 10500
$ws.Cells.Item(68, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(68, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(68, 16).Value = 420

# Row 69
$ws.Cells.Item(69, 4).Value = 44313
$ws.Cells.Item(69, 8).Value = 'Rodeo'
$ws.Cells.Item(69, 9).Value = '1a (cosecha)'
$ws.Cells.Item(69, 11).Value = 7500
$ws.Cells.Item(69, 12).Value = 8000
$ws.Cells.Item(69, 13).Value = 7750
$ws.Cells.Item(69, 16).Value = 310

# Row 70
$ws.Cells.Item(70, 4).Value = 44350
$ws.Cells.Item(70, 8).Value = 'Asterix'
$ws.Cells.Item(70, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(70, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(70, 15).Value = 'Región de La Araucanía'

# Row 71
$ws.Cells.Item(71, 4).Value = 44623
$ws.Cells.Item(71, 8).Value = 'Patagonia'
$ws.Cells.Item(71, 9).Value = '1a (cosecha)'
$ws.Cells.Item(71, 11).Value = 9000
$ws.Cells.Item(71, 12).Value = 10000
$ws.Cells.Item(71, 13).Value = 9500
$ws.Cells.Item(71, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(71, 16).Value = 380

# Row 72
$ws.Cells.Item(72, 4).Value = 44271
$ws.Cells.Item(72, 8).Value = 'Rodeo'
$ws.Cells.Item(72, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(72, 12).Value = 9500
$ws.Cells.Item(72, 13).Value = 9250
$ws.Cells.Item(72, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(72, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(72, 16).Value = 370

# --- Append new rows 73 and 74 ---

# Row 73
$ws.Cells.Item(73, 1).Value = 1
$ws.Cells.Item(73, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(73, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(73, 4).Value = 44273
$ws.Cells.Item(73, 5).Value = 15
$ws.Cells.Item(73, 6).Value = 100114001
$ws.Cells.Item(73, 7).Value = 'Papa'
$ws.Cells.Item(73, 8).Value = 'Rodeo'
$ws.Cells.Item(73, 9).Value = '1a (cosecha)'
$ws.Cells.Item(73, 10).Value = 1000
$ws.Cells.Item(73, 11).Value = 8000
$ws.Cells.Item(73, 12).Value = 8500
$ws.Cells.Item(73, 13).Value = 8250
$ws.Cells.Item(73, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(73, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(73, 16).Value = 330
$ws.Cells.Item(73, 17).Value = 25
$ws.Cells.Item(73, 18).Value = 'Hortaliza'
$ws.Cells.Item(73, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 74
$ws.Cells.Item(74, 1).Value = 1
$ws.Cells.Item(74, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(74, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(74, 4).Value = 44245
$ws.Cells.Item(74, 5).Value = 15
$ws.Cells.Item(74, 6).Value = 100114001
$ws.Cells.Item(74, 7).Value = 'Papa'
$ws.Cells.Item(74, 8).Value = 'Patagonia'
$ws.Cells.Item(74, 9).Value = '1a (cosecha)'
$ws.Cells.Item(74, 10).Value = 1000
$ws.Cells.Item(74, 11).Value = 8500
$ws.Cells.Item(74, 12).Value = 9000
$ws.Cells.Item(74, 13).Value = 8750
$ws.Cells.Item(74, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(74, 15).Value = 'Región de Los Lagos'
$ws.Cells.Item(74, 16).Value = 350
$ws.Cells.Item(74, 17).Value = 25
$ws.Cells.Item(74, 18).Value = 'Hortaliza'
$ws.Cells.Item(74, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

